$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Copy the style (date format) of A72 down to the two new rows
$ws.Range("A72").Copy() | Out-Null
$ws.Range("A73:A74").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 73: same task as row 72 ("Bugfixes")
$ws.Cells.Item(73, 1).Value = 45394
$ws.Cells.Item(73, 2).Value = 3
$ws.Cells.Item(73, 3).Value = "Bugfixes"

# Row 74: new task "Schreiben"
$ws.Cells.Item(74, 1).Value = 45394
$ws.Cells.Item(74, 2).Value = 3
$ws.Cells.Item(74, 3).Value = "Schreiben"

$ws.Range("C74").Select() | Out-Null

$wb.Save()
